$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing data rows (2..88) down by one row (to 3..89),
# working from the bottom up so we don't overwrite data before it's copied.
for ($r = 89; $r -ge 3; $r--) {
    $src = $r - 1
    for ($c = 1; $c -le 18; $c++) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($src, $c).Value2
    }
    # Column D carries a date-specific number format; keep it in sync with the
    # source row being copied down.
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($src, 4).NumberFormat
}

# Populate the newly freed row 2 with the new weekly record.
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "Vega Monumental Concepción"
$ws.Range("C2").Value = "Bíobío"
$ws.Range("D2").Value = 44631
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 100112024
$ws.Range("G2").Value = "Choclo"
$ws.Range("H2").Value = "Choclero"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 35000
$ws.Range("K2").Value = 150
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = 171
$ws.Range("N2").Value = "$/unidad"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 171
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"
